$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Apply number formats (reuse existing styles s="1" for dates, s="5" for times) ---
$ws.Range("A55").Copy()
$ws.Range("A62").PasteSpecial(-4122)

$ws.Range("G55").Copy()
$ws.Range("G62:G69").PasteSpecial(-4122)

$ws.Range("H55").Copy()
$ws.Range("H62:H69").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- Row 62 : 2017-03-22, Reviewer, "Har rettet Test Suite for OC15 " ---
$ws.Cells.Item(62, 1).Value = 42816
$ws.Cells.Item(62, 5).Value = "Reviewer"
$ws.Cells.Item(62, 6).Value = "Har rettet Test Suite for OC15 "
$ws.Cells.Item(62, 7).Value = 0.34375
$ws.Cells.Item(62, 8).Value = 0.36458333333333331

# --- Row 63 : Implementer, "Har implementeret OC12 design" ---
$ws.Cells.Item(63, 5).Value = "Implementer"
$ws.Cells.Item(63, 6).Value = "Har implementeret OC12 design"
$ws.Cells.Item(63, 7).Value = 0.35416666666666669
$ws.Cells.Item(63, 8).Value = 0.375

# --- Row 64 : Reviewer, "Har lavet review på OC13 kode" ---
$ws.Cells.Item(64, 5).Value = "Reviewer"
$ws.Cells.Item(64, 6).Value = "Har lavet review på OC13 kode"
$ws.Cells.Item(64, 7).Value = 0.3923611111111111
$ws.Cells.Item(64, 8).Value = 0.4201388888888889

# --- Row 65 : Implementer, "Har hjulpet med implementation af OC15+ OC16 design" ---
$ws.Cells.Item(65, 5).Value = "Implementer"
$ws.Cells.Item(65, 6).Value = "Har hjulpet med implementation af OC15+ OC16 design"
$ws.Cells.Item(65, 7).Value = 0.42708333333333331
$ws.Cells.Item(65, 8).Value = 0.47916666666666669

# --- Row 66 : "Har hjulpet med implementation af OC15 Test Suite " ---
$ws.Cells.Item(66, 6).Value = "Har hjulpet med implementation af OC15 Test Suite "
$ws.Cells.Item(66, 7).Value = 0.50347222222222221
$ws.Cells.Item(66, 8).Value = 0.53125

# --- Row 67 : "Har implementeret Test klasse for Inertimoment " ---
$ws.Cells.Item(67, 6).Value = "Har implementeret Test klasse for Inertimoment "
$ws.Cells.Item(67, 7).Value = 0.53472222222222221
$ws.Cells.Item(67, 8).Value = 0.59375

# --- Row 68 : "Har implementeret Test klasse for HalvProfilhoejde" ---
$ws.Cells.Item(68, 6).Value = "Har implementeret Test klasse for HalvProfilhoejde"
$ws.Cells.Item(68, 7).Value = 0.53472222222222221
$ws.Cells.Item(68, 8).Value = 0.59375

# --- Row 69 : "Har implementeret Test klasse for Boejningsmoment" ---
$ws.Cells.Item(69, 6).Value = "Har implementeret Test klasse for Boejningsmoment"
$ws.Cells.Item(69, 7).Value = 0.60069444444444442
$ws.Cells.Item(69, 8).Value = 0.62847222222222221

# --- Row 70 : total hours, moved down from the old row 61 ---
$ws.Cells.Item(70, 9).Value = 5.3

# --- Column F width widened slightly (closest reachable value; COM's
#     ColumnWidth setter on this host quantises to 1/6-character steps) ---
$ws.Columns.Item(6).ColumnWidth = 48.666666666666664

# --- Sheet view: scroll position + active selection follow the new last row ---
$excel.ActiveWindow.ScrollRow = 54
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("H70").Select()
